$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5121893
$ws.Range("J17").Value = 5441617.5
$ws.Range("L17").Value = 16324852.5
$ws.Range("N17").Value = -16325188.5
$ws.Range("H58").Value = 1484.1818
$ws.Range("J58").Value = 2308
$ws.Range("L58").Value = 6924
$ws.Range("N58").Value = -7224
$ws.Range("H126").Value = 90000
$ws.Range("J126").Value = 90000
$ws.Range("L126").Value = 90000
$ws.Range("N126").Value = -99880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9863.5
$ws.Range("I45").Value = 10253.533
$ws.Range("K45").Value = 10253.533
$ws.Range("M45").Value = -9876.532999999999
$ws.Range("H74").Value = 5962.032
$ws.Range("I74").Value = 6027.9165
$ws.Range("J74").Value = 5736.143
$ws.Range("K74").Value = 6027.9165
$ws.Range("L74").Value = 5736.143
$ws.Range("M74").Value = -5153.9165
$ws.Range("N74").Value = -7484.143
$ws.Range("H77").Value = 5962.032
$ws.Range("I77").Value = 6027.9165
$ws.Range("J77").Value = 5736.143
$ws.Range("K77").Value = 30139.5825
$ws.Range("L77").Value = 28680.715
$ws.Range("M77").Value = -25771.5825
$ws.Range("N77").Value = -37416.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 224.5
$ws.Range("I22").Value = 199.33333
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 199.33333
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -26.33332999999999
$ws.Range("N22").Value = -646
$ws.Range("H59").Value = 65000
$ws.Range("J59").Value = 65000
$ws.Range("L59").Value = 65000
$ws.Range("N59").Value = -66694
$ws.Range("H86").Value = 437873.7
$ws.Range("I86").Value = 716224.5
$ws.Range("J86").Value = 4883.5557
$ws.Range("K86").Value = 716224.5
$ws.Range("L86").Value = 4883.5557
$ws.Range("M86").Value = -715101.5
$ws.Range("N86").Value = -7129.5557
$ws.Range("H89").Value = 437873.7
$ws.Range("I89").Value = 716224.5
$ws.Range("J89").Value = 4883.5557
$ws.Range("K89").Value = 3581122.5
$ws.Range("L89").Value = 24417.7785
$ws.Range("M89").Value = -3575506.5
$ws.Range("N89").Value = -35649.7785

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1518.3043
$ws.Range("I5").Value = 575.3333
$ws.Range("J5").Value = 1659.75
$ws.Range("K5").Value = 1725.9999
$ws.Range("L5").Value = 4979.25
$ws.Range("M5").Value = -1613.9999
$ws.Range("N5").Value = -5203.25
$ws.Range("H38").Value = 1665.1904
$ws.Range("I38").Value = 287.88235
$ws.Range("J38").Value = 7518.75
$ws.Range("K38").Value = 863.6470499999999
$ws.Range("L38").Value = 22556.25
$ws.Range("M38").Value = -516.6470499999999
$ws.Range("N38").Value = -23250.25
$ws.Range("H40").Value = 62.761906
$ws.Range("I40").Value = 30.733334
$ws.Range("J40").Value = 142.83333
$ws.Range("K40").Value = 122.933336
$ws.Range("L40").Value = 571.33332
$ws.Range("M40").Value = -53.933336
$ws.Range("N40").Value = -709.33332
$ws.Range("H81").Value = 16690.31
$ws.Range("I81").Value = 10800
$ws.Range("J81").Value = 17917.459
$ws.Range("K81").Value = 32400
$ws.Range("L81").Value = 53752.37699999999
$ws.Range("M81").Value = -31277
$ws.Range("N81").Value = -55998.37699999999
$ws.Range("H84").Value = 16690.31
$ws.Range("I84").Value = 10800
$ws.Range("J84").Value = 17917.459
$ws.Range("K84").Value = 97200
$ws.Range("L84").Value = 161257.131
$ws.Range("M84").Value = -91584
$ws.Range("N84").Value = -172489.131
$ws.Range("H93").Value = 7419.3
$ws.Range("J93").Value = 7822
$ws.Range("L93").Value = 23466
$ws.Range("N93").Value = -27210
$ws.Range("H99").Value = 6699.5557
$ws.Range("J99").Value = 10299.5
$ws.Range("L99").Value = 30898.5
$ws.Range("N99").Value = -35390.5
$ws.Range("H100").Value = 9390.666999999999
$ws.Range("J100").Value = 10500
$ws.Range("L100").Value = 31500
$ws.Range("N100").Value = -33122
$ws.Range("H108").Value = 2199.4
$ws.Range("I108").Value = 1749.25
$ws.Range("K108").Value = 5247.75
$ws.Range("M108").Value = -2367.75
$ws.Range("H132").Value = 5884043
$ws.Range("I132").Value = 1241.8572
$ws.Range("K132").Value = 11176.7148
$ws.Range("M132").Value = -8646.7148
$ws.Range("H135").Value = 1518.3043
$ws.Range("I135").Value = 575.3333
$ws.Range("J135").Value = 1659.75
$ws.Range("K135").Value = 5177.9997
$ws.Range("L135").Value = 14937.75
$ws.Range("M135").Value = -2642.9997
$ws.Range("N135").Value = -20007.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 13330.5
$ws.Range("J96").Value = 13330.5
$ws.Range("L96").Value = 13330.5
$ws.Range("N96").Value = -18822.5
$ws.Range("H113").Value = 134845.53
$ws.Range("I113").Value = 223647.22
$ws.Range("K113").Value = 223647.22
$ws.Range("M113").Value = -221477.22
$ws.Range("H126").Value = 9410.625
$ws.Range("J126").Value = 17032
$ws.Range("L126").Value = 51096
$ws.Range("N126").Value = -56036

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 12434.714
$ws.Range("I33").Value = 12434.714
$ws.Range("K33").Value = 12434.714
$ws.Range("M33").Value = -12144.714
$ws.Range("H46").Value = 2114.5217
$ws.Range("I46").Value = 1628.5714
$ws.Range("K46").Value = 1628.5714
$ws.Range("M46").Value = -1440.5714
$ws.Range("H80").Value = 47050.6
$ws.Range("J80").Value = 47050.6
$ws.Range("L80").Value = 47050.6
$ws.Range("N80").Value = -49296.6
$ws.Range("H82").Value = 2097
$ws.Range("I82").Value = 1747.3334
$ws.Range("J82").Value = 2271.8333
$ws.Range("K82").Value = 1747.3334
$ws.Range("L82").Value = 2271.8333
$ws.Range("M82").Value = -1386.3334
$ws.Range("N82").Value = -2993.8333
$ws.Range("H83").Value = 47050.6
$ws.Range("J83").Value = 47050.6
$ws.Range("L83").Value = 141151.8
$ws.Range("N83").Value = -152383.8
$ws.Range("H85").Value = 2097
$ws.Range("I85").Value = 1747.3334
$ws.Range("J85").Value = 2271.8333
$ws.Range("K85").Value = 1747.3334
$ws.Range("L85").Value = 2271.8333
$ws.Range("M85").Value = -499.3334
$ws.Range("N85").Value = -4767.8333
$ws.Range("H87").Value = 2500000
$ws.Range("J87").Value = 2500000
$ws.Range("L87").Value = 2500000
$ws.Range("N87").Value = -2502246
$ws.Range("H90").Value = 2500000
$ws.Range("J90").Value = 2500000
$ws.Range("L90").Value = 7500000
$ws.Range("N90").Value = -7511232
$ws.Range("H100").Value = 2089
$ws.Range("I100").Value = 934
$ws.Range("J100").Value = 5554
$ws.Range("K100").Value = 934
$ws.Range("L100").Value = 5554
$ws.Range("M100").Value = -393
$ws.Range("N100").Value = -6636
$ws.Range("H123").Value = 39429
$ws.Range("J123").Value = 39429
$ws.Range("L123").Value = 39429
$ws.Range("N123").Value = -49229

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696
